# Day 25 class is added — append a new row to the schedule table with
# Date / Topics / YouTube link columns.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Add a brand new row at the bottom of the table.
$newRow = $t.Rows.Add()
$lastRow = $t.Rows.Count

# --- Column 1: Date ("27th June", with "th" superscripted) -----------
$cell1 = $t.Cell($lastRow, 1)
$cell1.Range.Text = "27th June"
$para1 = $cell1.Range.Paragraphs.Item(1)
$cellStart1 = $para1.Range.Start
$supRange = $d.Range($cellStart1 + 2, $cellStart1 + 4)
$supRange.Font.Superscript = $true

# --- Column 2: Topics -------------------------------------------------
$cell2 = $t.Cell($lastRow, 2)
$cell2.Range.Text = "HOF, map, filter, reduce"

# --- Column 3: YouTube hyperlink --------------------------------------
$cell3 = $t.Cell($lastRow, 3)
$videoUrl = "https://www.youtube.com/watch?v=Vd2Z2G3i_xE"
$cell3.Range.Text = $videoUrl
$cell3b = $t.Cell($lastRow, 3)
$para3 = $cell3b.Range.Paragraphs.Item(1)
$textRange = $d.Range($para3.Range.Start, $para3.Range.End - 1)
$null = $d.Hyperlinks.Add($textRange, $videoUrl)
